$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.403.90"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.55"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.42"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3723"
$ws.Range("E7").Value = "  -1.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.09"
$ws.Range("E8").Value = "  +0.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3603"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.237"
$ws.Range("E10").Value = "  -2.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08084"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.70"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.561"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001263"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.245"
$ws.Range("E16").Value = "  -2.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.627.72"
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.17"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06857"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.03"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.485"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "23.415.03"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.69"
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.419"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.995"
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.30"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.331"
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.67"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.290"
$ws.Range("E31").Value = "  -4.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.806.28"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.724"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9429"
$ws.Range("E34").Value = "  -2.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02805"
$ws.Range("E35").Value = "  +2.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.19"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2511"
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.07163"
$ws.Range("E38").Value = "  -5.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08748"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.021"
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.366"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6994"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.35"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.02"
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6462"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.311"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.001"
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07962"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.95"
$ws.Range("E50").Value = "  -3.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.189"
$ws.Range("E51").Value = "  -1.48%  "
